# Updated symbol list on Sun Dec 11 23:17:54 UTC 2022 with GitHub Actions
#
# This script reproduces the cell-level edits shown in the source diff.
# The "Price" column (D) holds numeric-looking values that are stored as
# text (inlineStr) in the workbook, so we prefix them with a leading
# apostrophe when assigning via .Value - this forces Excel to keep them
# as text instead of silently converting them to real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $ws.Range($Address).Value = "'" + $Text
}

# --- Price (column D) updates -------------------------------------------------
Set-TextValue "D2"  "285.77"
Set-TextValue "D3"  "21.25"
Set-TextValue "D4"  "6.452"
Set-TextValue "D5"  "0.06359"
Set-TextValue "D6"  "3.604"
Set-TextValue "D7"  "1.552"
Set-TextValue "D8"  "6.562"
Set-TextValue "D9"  "0.8207"
Set-TextValue "D10" "0.01415"
Set-TextValue "D11" "0.1679"
Set-TextValue "D12" "0.08669"
Set-TextValue "D13" "0.03669"
Set-TextValue "D14" "0.03206"
Set-TextValue "D15" "0.09203"
Set-TextValue "D16" "3.725"
Set-TextValue "D17" "0.001638"
Set-TextValue "D19" "0.006244"
Set-TextValue "D20" "0.006293"
Set-TextValue "D24" "2.271"
Set-TextValue "D25" "0.3356"
Set-TextValue "D26" "0.1261"
Set-TextValue "D40" "0.04771"
Set-TextValue "D41" "0.007139"

# --- Rows 42/43 swap: CEJI <-> BKEXToken --------------------------------------
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1118"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003454"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"

# --- Remaining price updates ---------------------------------------------------
Set-TextValue "D45" "0.00007094"
Set-TextValue "D48" "0.004522"

# --- Label change (Worst in 24h flag removed) -----------------------------------
$ws.Range("E49").Value = "48CryptobidCoinCBC"

Set-TextValue "D50" "0.01241"
